# Updated cryptos list on Sun Oct 15 09:51:45 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) values for each coin row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates. These look like numbers/localized numbers, so force
# the cell to stay plain text (NumberFormat "@") and then restore the default
# "Normal" style so no extra style index gets stamped onto the cell.
$priceUpdates = @{
    "D2"  = "27.053.95"
    "D3"  = "1.565.88"
    "D5"  = "208.62"
    "D8"  = "22.11"
    "D11" = "0.0858"
    "D12" = "1.561.97"
    "D13" = "3.78"
    "D15" = "27.046.95"
    "D16" = "61.90"
    "D17" = "0.0₃0706"
    "D19" = "215.96"
    "D22" = "9.19"
    "D24" = "153.93"
    "D26" = "15.05"
    "D32" = "3.20"
    "D33" = "1.423.31"
    "D35" = "1.61"
    "D36" = "2.34"
    "D37" = "0.0167"
    "D42" = "2.34"
    "D43" = "1.00"
    "D44" = "64.85"
    "D46" = "1.703.43"
    "D47" = "86.86"
    "D50" = "0.0961"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
    $cell.Style = "Normal"
}

# Volume(1h) column (E) updates - plain text percentages, safe to assign directly.
$volumeUpdates = @{
    "E2"  = "  +0.54%  "
    "E3"  = "  +1.06%  "
    "E4"  = "  +0.52%  "
    "E5"  = "  +1.13%  "
    "E6"  = "  +0.87%  "
    "E8"  = "  +0.08%  "
    "E9"  = "  +1.17%  "
    "E10" = "  +1.87%  "
    "E11" = "  +0.40%  "
    "E12" = "  +0.80%  "
    "E13" = "  +1.33%  "
    "E14" = "  +0.39%  "
    "E15" = "  +0.54%  "
    "E16" = "  +0.45%  "
    "E17" = "  +1.37%  "
    "E18" = "  +2.51%  "
    "E21" = "  +2.41%  "
    "E22" = "  -0.34%  "
    "E23" = "  -0.27%  "
    "E24" = "  -0.22%  "
    "E25" = "  -0.21%  "
    "E26" = "  +0.80%  "
    "E27" = "  +1.62%  "
    "E28" = "  +0.47%  "
    "E29" = "  +1.67%  "
    "E30" = "  +4.34%  "
    "E31" = "  +0.69%  "
    "E32" = "  +3.60%  "
    "E33" = "  +0.22%  "
    "E34" = "  +12.74%  "
    "E35" = "  +1.75%  "
    "E36" = "  +2.76%  "
    "E37" = "  +1.39%  "
    "E38" = "  +1.82%  "
    "E39" = "  +2.16%  "
    "E40" = "  +0.41%  "
    "E41" = "  +0.39%  "
    "E42" = "  +0.60%  "
    "E43" = "  +0.47%  "
    "E44" = "  +0.52%  "
    "E45" = "  -0.49%  "
    "E46" = "  +1.14%  "
    "E47" = "  -0.63%  "
    "E48" = "  +3.42%  "
    "E49" = "  +0.98%  "
    "E50" = "  +0.10%  "
    "E51" = "  +0.46%  "
}

foreach ($cellRef in $volumeUpdates.Keys) {
    $ws.Range($cellRef).Value = $volumeUpdates[$cellRef]
}
